# Adds a new "Testing for Max Length" test-case section to the
# PredictionTestCases sheet (rows 73-77), following the same layout
# pattern used by the other test-case sections already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 73: section title, styled like the other section headers (e.g. A63) ---
$ws.Range("A73").Value = "Testing for  Max Length"
$ws.Range("A63").Copy()
$ws.Range("A73").PasteSpecial(-4122)   # xlPasteFormats - copy formatting only

# --- Row 74: column headers "max length " / "Result", bold (no fill) ---
$ws.Range("B74").Value = "max length "
$ws.Range("C74").Value = "Result"
$ws.Range("B74:C74").Font.Bold = $true

# --- Rows 75-77: test case data ---
# Values are entered in the specific order below so that the shared string
# table ends up populated in the same sequence as the source workbook.
$ws.Range("A75").Value = "Credit Score"
$ws.Range("A76").Value = "Origination UPB"
$ws.Range("C75").Value = "UI doesn't allow User to enter"
$ws.Range("A77").Value = "Loan Term"

$ws.Range("B75").Value = 3
$ws.Range("B76").Value = 9
$ws.Range("B77").Value = 3

$ws.Range("C76").Value = "UI doesn't allow User to enter"
$ws.Range("C77").Value = "UI doesn't allow User to enter"

# Update the active selection, matching the saved workbook state.
$ws.Range("A9").Select()
